$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a number (e.g. "1.000", "27.276.09")
# must be force-typed as Text so Excel does not silently coerce them to
# a Double (dropping trailing zeros) or a date. We flip NumberFormat to
# "@" for the assignment, then restore the cell to the default "Normal"
# style so no stray style index is left behind on the cell.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.276.09'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.820.02'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.14%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.34'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4671'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.51%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3775'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.44%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07409'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.84%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8719'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.31%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.63'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.40%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.821.50'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.18%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.691'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.414'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.53%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.50'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.07084'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.000'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008767'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.93%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('E20').Value = '  +0.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.284.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.312'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.95%  '
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.050.96'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.70%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.938'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.43'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.29%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.246'
$ws.Range('D27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '18.54'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.316'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '117.08'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08944'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.83%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7839'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.70%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.183'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.530'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.60%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.915'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.000'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.06%  '
$ws.Range('E37').Value = '  +1.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01971'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05253'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.280'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.34%  '
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.382'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +20.99%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5321'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.58%  '
$ws.Range('E43').Value = '  +1.81%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.1693'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.49%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.619'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.03%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5065'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.52%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.43'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '105.65'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.53%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.9998'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.671'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.41%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06331'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.64%  '
